# Weekly fruit/vegetable data refresh:
# - A new weekly price reading (date 2022-06-07 / serial 44719) is recorded for
#   the first "Primera"/"Segunda" pair (rows 114-115), whose Origen is updated
#   to "Región Metropolitana".
# - All subsequent existing pairs (rows 114-198) shift down by one pair (2 rows),
#   so the old content of rows 114-198 ends up, unchanged, in rows 116-200.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Snapshot the existing block of rows (114-198) before touching anything.
$src = $ws.Range("A114:R198")
$arr = $src.Value()

# 2. Make sure the destination's Fecha (date) column already carries the same
#    number format as the rest of the column before values land on it -
#    otherwise the brand-new rows (199-200) would pick up a default date
#    format instead of matching the existing column formatting.
$ws.Range("D116:D200").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# 3. Shift that whole block down by 2 rows (to 116-200).
$dest = $ws.Range("A116:R200")
$dest.Value = $arr

# 4. Overwrite the freed-up first pair (114-115) with the new weekly reading:
#    only the date and the origin region change; the rest of the row keeps
#    the values it already had.
$newDate = Get-Date -Year 2022 -Month 6 -Day 7 -Hour 0 -Minute 0 -Second 0
$ws.Range("D114").Value = $newDate
$ws.Range("D115").Value = $newDate
$ws.Range("O114").Value = "Región Metropolitana"
$ws.Range("O115").Value = "Región Metropolitana"
